$d = $word.ActiveDocument

# Anchor on the "git stash clear" paragraph (paragraph 21) and insert a new
# paragraph after it; subsequent paragraphs are appended the same way so the
# five new bullet lines land between "git stash clear" and the following
# blank paragraph.
$anchor = $d.Paragraphs.Item(21)
$anchor.Range.InsertParagraphAfter()

$texts = @(
  "git remote add origin “url of repo” – Adding a url named origin on our folder",
  "git remote -v – List out url of the remote repositories ",
  "git push [urlName] [branchName] – Push the commits to a repo with selected branch name and url name",
  "git branch – List out all the branches in current git repo",
  "git branch [branch Name] – Create a new branch [branch Name]"
)

$startIndex = 22
for ($i = 0; $i -lt $texts.Count; $i++) {
    $p = $d.Paragraphs.Item($startIndex + $i)
    $p.Range.Text = $texts[$i]
    if ($i -lt $texts.Count - 1) {
        $p2 = $d.Paragraphs.Item($startIndex + $i)
        $p2.Range.InsertParagraphAfter()
    }
}
